{"js": "// Insert a new list item \"Clean up repo and make Tree standalone (and maybe\n// other potential libraries)\" right after the \"Implement the simplest\n// solver\" list item (same numbered list, numId=7).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.trim() === \"Implement the simplest solver\") {\n    anchor = p;\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not locate paragraph \"Implement the simplest solver\"');\n}\n\nconst newPara = anchor.insertParagraph(\n  \"Clean up repo and make Tree standalone (and maybe other potential libraries)\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Insert a new list item \"Clean up repo and make Tree standalone (and maybe\n# other potential libraries)\" right after the \"Implement the simplest\n# solver\" list item (same numbered list, numId=7).\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r`a\", \"`r\", \"`n\") -eq \"Implement the simplest solver\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate paragraph 'Implement the simplest solver'\"\n}\n\n$newRange = $target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = \"Clean up repo and make Tree standalone (and maybe other potential libraries)\"\n"}
